# Updates cryptos list values (Price / Volume(1h)) to match the latest
# scrape, and swaps the ordi / BitcoinSV rows (47-48) per the new ranking.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "44.824.78"
$ws.Range("E2").Value = "  +0.05%  "
$ws.Range("D3").Value = "2.261.85"
$ws.Range("E3").Value = "  +0.67%  "
$ws.Range("E4").Value = "  -0.75%  "
$ws.Range("D5").Value = "'301.96"
$ws.Range("E5").Value = "  -1.54%  "
$ws.Range("D6").Value = "'94.05"
$ws.Range("E6").Value = "  -1.62%  "
$ws.Range("E7").Value = "  -1.55%  "
$ws.Range("D8").Value = "'1.00"
$ws.Range("E8").Value = "  -0.56%  "
$ws.Range("D9").Value = "'0.509"
$ws.Range("E9").Value = "  -2.25%  "
$ws.Range("D10").Value = "'34.16"
$ws.Range("E10").Value = "  -2.80%  "
$ws.Range("D11").Value = "'0.0789"
$ws.Range("E11").Value = "  -1.97%  "
$ws.Range("E12").Value = "  -0.42%  "
$ws.Range("E13").Value = "  -1.33%  "
$ws.Range("D14").Value = "2.607.10"
$ws.Range("E14").Value = "  +0.59%  "
$ws.Range("D15").Value = "2.259.51"
$ws.Range("E15").Value = "  -3.01%  "
$ws.Range("D16").Value = "'13.57"
$ws.Range("E16").Value = "  -0.29%  "
$ws.Range("D17").Value = "'0.799"
$ws.Range("E17").Value = "  -4.91%  "
$ws.Range("D18").Value = "44.757.12"
$ws.Range("E18").Value = "  +0.35%  "
$ws.Range("D19").Value = "'12.90"
$ws.Range("E19").Value = "  +7.76%  "
$ws.Range("D20").Value = "0.0₃0924"
$ws.Range("E20").Value = "  -2.70%  "
$ws.Range("E21").Value = "  -3.57%  "
$ws.Range("D22").Value = "'65.62"
$ws.Range("E22").Value = "  +0.25%  "
$ws.Range("D23").Value = "'238.04"
$ws.Range("E23").Value = "  -0.65%  "
$ws.Range("E24").Value = "  -2.78%  "
$ws.Range("E25").Value = "  -0.47%  "
$ws.Range("E26").Value = "  -4.86%  "
$ws.Range("D27").Value = "'41.31"
$ws.Range("E27").Value = "  +9.22%  "
$ws.Range("D28").Value = "'2.30"
$ws.Range("E28").Value = "  -0.24%  "
$ws.Range("D29").Value = "'9.61"
$ws.Range("E29").Value = "  -2.57%  "
$ws.Range("D30").Value = "'19.57"
$ws.Range("E30").Value = "  -2.06%  "
$ws.Range("D31").Value = "'152.16"
$ws.Range("E31").Value = "  +0.94%  "
$ws.Range("D32").Value = "'5.54"
$ws.Range("E32").Value = "  -8.31%  "
$ws.Range("D33").Value = "'0.0790"
$ws.Range("E33").Value = "  -0.83%  "
$ws.Range("D34").Value = "'2.55"
$ws.Range("E34").Value = "  -3.07%  "
$ws.Range("E35").Value = "  -3.49%  "
$ws.Range("E36").Value = "  -1.44%  "
$ws.Range("E37").Value = "  -3.47%  "
$ws.Range("D38").Value = "'1.76"
$ws.Range("E38").Value = "  -5.04%  "
$ws.Range("D39").Value = "'3.94"
$ws.Range("E39").Value = "  +3.94%  "
$ws.Range("E40").Value = "  +2.58%  "
$ws.Range("D41").Value = "'3.21"
$ws.Range("E41").Value = "  -5.27%  "
$ws.Range("D42").Value = "'13.54"
$ws.Range("E42").Value = "  -10.57%  "
$ws.Range("D43").Value = "'0.999"
$ws.Range("E43").Value = "  -0.87%  "
$ws.Range("D44").Value = "'1.91"
$ws.Range("E44").Value = "  +8.60%  "
$ws.Range("D45").Value = "1.732.41"
$ws.Range("E45").Value = "  -6.36%  "
$ws.Range("E46").Value = "  +2.11%  "
$ws.Range("B47").Value = "BitcoinSV"
$ws.Range("C47").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D47").Value = "'75.93"
$ws.Range("E47").Value = "  -4.77%  "
$ws.Range("B48").Value = "ordi"
$ws.Range("C48").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range("D48").Value = "'69.25"
$ws.Range("E48").Value = "  +0.13%  "
$ws.Range("D49").Value = "'95.58"
$ws.Range("E49").Value = "  -3.39%  "
$ws.Range("D50").Value = "'53.41"
$ws.Range("E50").Value = "  -1.93%  "
$ws.Range("E51").Value = "  -4.52%  "
